$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.674.38"
$ws.Range("E2").Value = "  +3.75%  "
$ws.Range("D3").Value = "1.915.85"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'250.69"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'0.702"
$ws.Range("E6").Value = "  +2.92%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'44.90"
$ws.Range("E8").Value = "  +2.91%  "
$ws.Range("D9").Value = "'0.372"
$ws.Range("E9").Value = "  +4.27%  "
$ws.Range("D10").Value = "'58.44"
$ws.Range("E10").Value = "  +9.61%  "
$ws.Range("D11").Value = "'0.0763"
$ws.Range("E11").Value = "  +2.59%  "
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "'14.55"
$ws.Range("E13").Value = "  +7.74%  "
$ws.Range("E14").Value = "  +6.91%  "
$ws.Range("D15").Value = "2.196.17"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "'5.13"
$ws.Range("E16").Value = "  +4.05%  "
$ws.Range("D17").Value = "1.922.33"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "36.664.80"
$ws.Range("E18").Value = "  +3.71%  "
$ws.Range("D19").Value = "'74.69"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +4.63%  "
$ws.Range("D21").Value = "'250.49"
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").Value = "'13.39"
$ws.Range("E22").Value = "  +4.46%  "
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("D24").Value = "'2.64"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("D27").Value = "'168.75"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").Value = "'8.78"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "'18.73"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  +6.58%  "
$ws.Range("D32").Value = "'0.0620"
$ws.Range("E32").Value = "  +4.58%  "
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "'0.0900"
$ws.Range("E34").Value = "  +23.35%  "
$ws.Range("D35").Value = "'1.91"
$ws.Range("E35").Value = "  +6.94%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +4.83%  "
$ws.Range("D38").Value = "'0.879"
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("D39").Value = "'17.68"
$ws.Range("E39").Value = "  +48.95%  "
$ws.Range("D40").Value = "'2.02"
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("D41").Value = "'106.11"
$ws.Range("E41").Value = "  +9.60%  "
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("D43").Value = "'17.52"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "'2.93"
$ws.Range("E44").Value = "  +22.40%  "
$ws.Range("D45").Value = "'1.10"
$ws.Range("E45").Value = "  +3.42%  "
$ws.Range("D46").Value = "1.345.00"
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").Value = "'2.38"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'0.0813"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").Value = "'2.80"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'6.44"
$ws.Range("E50").Value = "  +2.24%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'43.78"
$ws.Range("E51").Value = "  +4.16%  "
